$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.08461582660675
$ws.Range("B1").Value = 2.143845558166504
$ws.Range("C1").Value = 9.339086532592773
$ws.Range("D1").Value = 1.009033441543579
$ws.Range("E1").Value = 1.032935738563538
